$wb = $excel.ActiveWorkbook

$qaSheet = $wb.Worksheets.Item("qa")
$automationSheet = $wb.Worksheets.Item("automation")

# Update the "qa" sheet's login URL (B4) to match the current automation login URL.
# This leaves the older "bidline" URL with no remaining references in the shared
# string table, so it is dropped on save.
$qaSheet.Activate() | Out-Null
$qaSheet.Range("B4").Value = "https://crewbid-automation.firebaseapp.com/login"
$qaSheet.Range("B4").Select() | Out-Null

# Restore "automation" as the active sheet (it was active before this edit), with
# its B4 cell still selected.
$automationSheet.Activate() | Out-Null
$automationSheet.Range("B4").Select() | Out-Null
